$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-free inline approach: for numeric-looking Price values we must force
# the cell to Text format first so Excel does not coerce "0.9995" etc. into a
# number, then reset the style back to Normal so no stray formatting remains.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '30.714.98'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.44%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.890.99'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +1.71%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9995'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.03%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '237.57'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.60%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.9995'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.01%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4847'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +1.16%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.2870'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +2.25%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.06562'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +1.28%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.900.51'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.30%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07447'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.86%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '16.63'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.09%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.099'
$ws.Range("D13").Style = "Normal"

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '88.03'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.92%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.6646'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +2.89%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '30.658.47'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.43%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '13.22'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +0.52%  '

$ws.Range("E18").Value = '  -0.05%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.000007618'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.69%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '230.47'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +2.78%  '

$ws.Range("B21").Value = 'BinanceUSD'
$ws.Range("C21").Value = 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.000'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.00%  '

$ws.Range("B22").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C22").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '2.078.46'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.75%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '5.277'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.20%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '6.210'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.35%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.432'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +2.62%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '168.00'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.71%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.72'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.53%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.956'
$ws.Range("D28").Style = "Normal"

$ws.Range("E29").Value = '  +10.94%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.391'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -3.00%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '4.337'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +1.96%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.028'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.70%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.05055'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.83%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.207'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +5.58%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.7542'
$ws.Range("D35").Style = "Normal"

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.9992'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.16%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.710'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.92%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01898'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +3.50%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.665'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +2.64%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.9189'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.43%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.062'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.05%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '107.19'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +1.08%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.4283'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.82%  '

$ws.Range("E44").Value = '  +0.17%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '5.653'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.27%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.446'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.29%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '64.86'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +2.53%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.1273'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -1.62%  '

$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.011'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +3.50%  '

$ws.Range("B50").Value = 'NEARProtocol'
$ws.Range("C50").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.481'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -0.43%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '34.09'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.91%  '
